$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1) Column widths (A:H explicit widths, J:M shared width) - mirrors the new
#    <cols> block. ColumnWidth takes "characters"; the engine re-serialises
#    it into xlsx "width" units with a fixed +5/6 offset, so we back the
#    character value out of the desired xlsx width.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth  = 3.3307291666666665   # -> width 4.1640625
$ws.Columns.Item(2).ColumnWidth  = 2.9986979166666665   # -> width 3.83203125
$ws.Columns.Item(3).ColumnWidth  = 7.998697916666667    # -> width 8.83203125
$ws.Columns.Item(4).ColumnWidth  = 8.998697916666666    # -> width 9.83203125
$ws.Columns.Item(5).ColumnWidth  = 13.666666666666666   # -> width 14.5
$ws.Columns.Item(6).ColumnWidth  = 11.998697916666666   # -> width 12.83203125
$ws.Columns.Item(7).ColumnWidth  = 11.830729166666666   # -> width 12.6640625
$ws.Columns.Item(8).ColumnWidth  = 10.666666666666666   # -> width 11.5
$ws.Columns.Item(10).ColumnWidth = 2.9986979166666665   # -> width 3.83203125
$ws.Columns.Item(11).ColumnWidth = 2.9986979166666665   # -> width 3.83203125
$ws.Columns.Item(12).ColumnWidth = 2.9986979166666665   # -> width 3.83203125
$ws.Columns.Item(13).ColumnWidth = 2.9986979166666665   # -> width 3.83203125

# ---------------------------------------------------------------------------
# 2) Re-format a handful of cells in the plate-layout block (E:H) whose fill
#    formatting shifts between the "bold/coloured" and "plain" variants used
#    elsewhere in the same columns. Borrow the look from a sibling cell that
#    already carries the wanted format so the xlsx reuses the existing style
#    record instead of inventing a new one.
# ---------------------------------------------------------------------------
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)

$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial(-4122)

$ws.Range("D5").Copy()
$ws.Range("F5").PasteSpecial(-4122)

$ws.Range("D6").Copy()
$ws.Range("E6").PasteSpecial(-4122)

$ws.Range("D7").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("G7").PasteSpecial(-4122)

$ws.Range("D8").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("D8").Copy()
$ws.Range("G8").PasteSpecial(-4122)

$ws.Range("D9").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("I9").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("H9").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Relabel the plate-layout condition columns: "lacto" qualifier added to
#    the anaerobic/aerobic/O2-limit/fodder-yeast conditions (plus their typos
#    fixed along the way), and the stray "AABDH"/"DMR" columns collapsed onto
#    the existing "H20" label per the notes-sheet clarification.
#
#    Cell-write order below matches how the shared-string table grows (new
#    unique strings 43..50 get interned in this exact sequence): E/F/G for
#    rows 2-5, then E/F/G for rows 6-9, then H for rows 6-9, then H for rows
#    2-5 last.
# ---------------------------------------------------------------------------
foreach ($r in 2,3,4,5) {
    $ws.Range("E$r").Value = " Anaerobic lacto"
    $ws.Range("F$r").Value = " Aerobic lacto"
    $ws.Range("G$r").Value = "O2 limit lacto"
}
foreach ($r in 6,7,8,9) {
    $ws.Range("E$r").Value = " Anaerobic lacto_Blk"
    $ws.Range("F$r").Value = " Aerobic lacto_Blk"
    $ws.Range("G$r").Value = "O2 limit lacto_Blk"
}
foreach ($r in 6,7,8,9) {
    $ws.Range("H$r").Value = "Fodder yeast_Blk"
}
foreach ($r in 2,3,4,5) {
    $ws.Range("H$r").Value = "Fodder yeast"
}
foreach ($r in 2,3,4,5,6,7,8,9) {
    $ws.Range("J$r").Value = "H20"
    $ws.Range("K$r").Value = "H20"
}

# ---------------------------------------------------------------------------
# 4) Selection moved from the old K2:K9 block down to H10.
# ---------------------------------------------------------------------------
$ws.Range("H10").Select()
